# Applies the "output generated at c8c62b6" update:
#  - bumps the header date
#  - replaces each three-digit x one-digit multiplication prompt with a new one

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

# Header date line
Replace-Text "2025-07-04 Friday" "2025-07-05 Saturday"

# Row 1 (problems 1-5)
Replace-Text "517×4=" "782×8="
Replace-Text "986×5=" "542×4="
Replace-Text "357×5=" "775×9="
Replace-Text "465×7=" "988×8="
Replace-Text "372×4=" "543×2="

# Row 2 (problems 6-10)
Replace-Text "687×4=" "292×9="
Replace-Text "716×7=" "489×6="
Replace-Text "551×7=" "475×9="
Replace-Text "272×3=" "935×9="
Replace-Text "903×7=" "605×7="

# Row 3 (problems 11-15)
Replace-Text "557×2=" "861×9="
Replace-Text "332×7=" "538×4="
Replace-Text "891×2=" "627×2="
Replace-Text "715×6=" "472×5="
Replace-Text "514×7=" "919×4="

# Row 4 (problems 16-20) -- first "683×5=" (column 3) becomes "509×3="
Replace-Text "223×7=" "849×7="
Replace-Text "494×6=" "278×7="
$t = $d.Tables.Item(1)
$t.Cell(15, 3).Range.Text = "509×3="
Replace-Text "393×7=" "610×3="
Replace-Text "685×6=" "178×7="

# Row 5 (problems 21-25) -- second "683×5=" (column 1) becomes "536×9="
$t.Cell(20, 1).Range.Text = "536×9="
Replace-Text "986×3=" "560×6="
Replace-Text "319×4=" "636×4="
Replace-Text "464×5=" "266×2="
Replace-Text "967×3=" "656×9="
